$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "33.964.18"
$ws.Cells.Item(2, 5).Value = "  -0.26%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "1.772.97"
$ws.Cells.Item(3, 5).Value = "  -2.15%  "

# Row 4
$ws.Cells.Item(4, 5).Value = "  +0.24%  "

# Row 5
$ws.Cells.Item(5, 4).Value = "'225.20"
$ws.Cells.Item(5, 5).Value = "  -1.22%  "

# Row 6
$ws.Cells.Item(6, 4).Value = "'0.547"

# Row 7
$ws.Cells.Item(7, 5).Value = "  +0.25%  "

# Row 8
$ws.Cells.Item(8, 4).Value = "'31.22"
$ws.Cells.Item(8, 5).Value = "  -0.44%  "

# Row 9
$ws.Cells.Item(9, 4).Value = "'0.279"
$ws.Cells.Item(9, 5).Value = "  -0.59%  "

# Row 10
$ws.Cells.Item(10, 4).Value = "'0.0654"
$ws.Cells.Item(10, 5).Value = "  -1.76%  "

# Row 11
$ws.Cells.Item(11, 4).Value = "'0.0927"
$ws.Cells.Item(11, 5).Value = "  +0.03%  "

# Row 12
$ws.Cells.Item(12, 4).Value = "2.025.91"
$ws.Cells.Item(12, 5).Value = "  -2.19%  "

# Row 13
$ws.Cells.Item(13, 5).Value = "  +7.05%  "

# Row 14
$ws.Cells.Item(14, 4).Value = "1.770.74"
$ws.Cells.Item(14, 5).Value = "  -2.15%  "

# Row 15
$ws.Cells.Item(15, 4).Value = "33.955.15"
$ws.Cells.Item(15, 5).Value = "  -0.14%  "

# Row 16
$ws.Cells.Item(16, 4).Value = "'0.620"
$ws.Cells.Item(16, 5).Value = "  -3.01%  "

# Row 17
$ws.Cells.Item(17, 4).Value = "'4.19"
$ws.Cells.Item(17, 5).Value = "  -1.54%  "

# Row 18
$ws.Cells.Item(18, 4).Value = "'68.51"
$ws.Cells.Item(18, 5).Value = "  -1.15%  "

# Row 19
$ws.Cells.Item(19, 4).Value = "'251.38"
$ws.Cells.Item(19, 5).Value = "  -2.28%  "

# Row 20
$ws.Cells.Item(20, 4).Value = "0.0₃0734"
$ws.Cells.Item(20, 5).Value = "  -1.57%  "

# Row 21
$ws.Cells.Item(21, 5).Value = "  +0.30%  "

# Row 22
$ws.Cells.Item(22, 4).Value = "'10.27"
$ws.Cells.Item(22, 5).Value = "  -2.04%  "

# Row 23
$ws.Cells.Item(23, 5).Value = "  -3.53%  "

# Row 24
$ws.Cells.Item(24, 5).Value = "  -2.94%  "

# Row 25
$ws.Cells.Item(25, 4).Value = "'156.05"
$ws.Cells.Item(25, 5).Value = "  -1.40%  "

# Row 26
$ws.Cells.Item(26, 4).Value = "'16.33"
$ws.Cells.Item(26, 5).Value = "  -1.46%  "

# Row 27
$ws.Cells.Item(27, 4).Value = "'6.97"
$ws.Cells.Item(27, 5).Value = "  -1.92%  "

# Row 28
$ws.Cells.Item(28, 5).Value = "  -1.45%  "

# Row 29
$ws.Cells.Item(29, 5).Value = "  +0.28%  "

# Row 30
$ws.Cells.Item(30, 4).Value = "'3.75"
$ws.Cells.Item(30, 5).Value = "  -3.06%  "

# Row 31
$ws.Cells.Item(31, 4).Value = "'0.0508"
$ws.Cells.Item(31, 5).Value = "  -0.61%  "

# Row 32
$ws.Cells.Item(32, 5).Value = "  -0.88%  "

# Row 33
$ws.Cells.Item(33, 4).Value = "'3.55"
$ws.Cells.Item(33, 5).Value = "  +1.32%  "

# Row 34
$ws.Cells.Item(34, 4).Value = "'1.83"
$ws.Cells.Item(34, 5).Value = "  +2.13%  "

# Row 35
$ws.Cells.Item(35, 4).Value = "1.445.14"
$ws.Cells.Item(35, 5).Value = "  -6.00%  "

# Row 36
$ws.Cells.Item(36, 5).Value = "  -1.88%  "

# Row 37
$ws.Cells.Item(37, 4).Value = "'0.622"
$ws.Cells.Item(37, 5).Value = "  +0.12%  "

# Row 38
$ws.Cells.Item(38, 5).Value = "  -0.84%  "

# Row 39
$ws.Cells.Item(39, 4).Value = "'2.84"
$ws.Cells.Item(39, 5).Value = "  +1.44%  "

# Row 40
$ws.Cells.Item(40, 4).Value = "'82.39"
$ws.Cells.Item(40, 5).Value = "  -2.31%  "

# Row 41
$ws.Cells.Item(41, 5).Value = "  +0.41%  "

# Row 42
$ws.Cells.Item(42, 4).Value = "'0.882"
$ws.Cells.Item(42, 5).Value = "  -2.71%  "

# Row 43
$ws.Cells.Item(43, 5).Value = "  -4.72%  "

# Row 44
$ws.Cells.Item(44, 5).Value = "  -2.48%  "

# Row 45
$ws.Cells.Item(45, 5).Value = "  -1.89%  "

# Row 46
$ws.Cells.Item(46, 4).Value = "1.926.15"
$ws.Cells.Item(46, 5).Value = "  -2.20%  "

# Row 47
$ws.Cells.Item(47, 2).Value = "FraxShare"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Cells.Item(47, 4).Value = "'5.72"
$ws.Cells.Item(47, 5).Value = "  -0.07%  "

# Row 48
$ws.Cells.Item(48, 2).Value = "PaxDollar"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Cells.Item(48, 4).Value = "'1.00"
$ws.Cells.Item(48, 5).Value = "  +0.37%  "

# Row 49
$ws.Cells.Item(49, 4).Value = "'11.83"
$ws.Cells.Item(49, 5).Value = "  +2.55%  "

# Row 50
$ws.Cells.Item(50, 4).Value = "'97.05"
$ws.Cells.Item(50, 5).Value = "  +2.70%  "

# Row 51
$ws.Cells.Item(51, 4).Value = "'49.35"
$ws.Cells.Item(51, 5).Value = "  -5.80%  "
